$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new issue log entry in column A, row 12
$ws.Range("A12").Value = "dynamically change position with event listener "

# Reflect the trailing selection change left in the saved file
$ws.Range("A15").Select()
